$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 13 ("2021年") below the existing last data row (12).
# Copy row 12 first so formatting/style and the "blank" cells (Q, R, T, Z,
# AC, AD, AG, AH, AS, BK) carry over exactly, then overwrite with the new
# year's values.
$ws.Range("A12:BS12").Copy($ws.Range("A13:BS13"))

$ws.Range("A13").Value = "2021年"

$values = @{
    "B"  = 159346
    "C"  = 11375
    "D"  = 32443
    "E"  = 15955
    "F"  = 26042
    "G"  = 2702
    "H"  = 470
    "I"  = 47103
    "J"  = 2115
    "K"  = 22582
    "L"  = 589
    "M"  = 10523
    "N"  = 214
    "O"  = 52
    "P"  = 9710
    "S"  = 582
    "U"  = 19
    "V"  = 39
    "W"  = 5683
    "X"  = 451
    "Y"  = 234
    "AA" = 796
    "AB" = 553
    "AE" = 1144
    "AF" = 126
    "AI" = 598
    "AJ" = 1781
    "AK" = 50031
    "AL" = 6462
    "AM" = 7994
    "AN" = 12910
    "AO" = 7653
    "AP" = 791
    "AQ" = 391
    "AR" = 9185
    "AT" = 232
    "AU" = 162
    "AV" = 465
    "AW" = 440
    "AX" = 2331
    "AY" = 19037
    "AZ" = 5416
    "BA" = 56529
    "BB" = 2347
    "BC" = 1633
    "BD" = 2456
    "BE" = 7070
    "BF" = 1144
    "BG" = 27
    "BH" = 142
    "BI" = 4272
    "BJ" = 932
    "BL" = 4841
    "BM" = 1661
    "BN" = 29853
    "BO" = 1290
    "BP" = 744
    "BQ" = 5306
    "BR" = 2718
    "BS" = 34004
}

foreach ($col in $values.Keys) {
    $ws.Range("$col" + "13").Value = $values[$col]
}
